$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted as row 29 (pushing every
# existing row from 29 down through 138 down by one, to 30-139).
$ws.Rows("29:29").Insert()

# Populate the newly inserted row 29 with the new observation. All the
# descriptive columns match the neighbouring rows for this
# market/product (Vega Monumental Concepción - Ajo, Chino, Primera).
$ws.Cells.Item(29, 1).Value = 11
$ws.Cells.Item(29, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(29, 3).Value = "Bíobío"
$ws.Cells.Item(29, 4).Value = 44600
$ws.Cells.Item(29, 5).Value = 8
$ws.Cells.Item(29, 6).Value = 100112003
$ws.Cells.Item(29, 7).Value = "Ajo"
$ws.Cells.Item(29, 8).Value = "Chino"
$ws.Cells.Item(29, 9).Value = "Primera"
$ws.Cells.Item(29, 10).Value = 250
$ws.Cells.Item(29, 11).Value = 14000
$ws.Cells.Item(29, 12).Value = 15000
$ws.Cells.Item(29, 13).Value = 14520
$ws.Cells.Item(29, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(29, 15).Value = "China"
$ws.Cells.Item(29, 16).Value = 1452
$ws.Cells.Item(29, 17).Value = 10
$ws.Cells.Item(29, 18).Value = "Hortaliza"
